# Shop.xlsx edit: re-insert a "key|Int" column (G) into ShopProductTable,
# shifting the existing reward columns right by one, and populate the new
# trailing "테이블연결" / "Jason화" helper columns (W/X) with the running
# JSON blob used by the lookup table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ShopProductTable")

# --- 1. Insert the new "key|Int" column at G, pushing rewardType1.. right ---
$ws.Columns("G:G").Insert()
$ws.Columns("G:G").ColumnWidth = 5.8

$ws.Range("G1").Value = "key|Int"

$keys = @{
    2 = 434
    3 = 806
    4 = 548
    5 = 537
    6 = 314
    7 = 876
    8 = 973
    9 = 180
    10 = 721
    11 = 884
    12 = 217
    13 = 394
    14 = 612
    15 = 501
    16 = 930
    17 = 959
}
foreach ($r in $keys.Keys) {
    $ws.Cells.Item($r, 7).Value = $keys[$r]
}

# --- 2. New trailing columns W (테이블연결) / X (Jason화) ---
$ws.Range("W1").Value = "테이블연결"
$ws.Range("X1").Value = "Jason화"

$wVals = @{
    5 = ',{"id":"ev5_oneplustwo_2","key":537,"tp1":"cu","vl1":"EN","cn1":300}'
    6 = ',{"id":"ev5_oneplustwo_2","key":537,"tp1":"cu","vl1":"EN","cn1":300},{"id":"ev5_oneplustwo_3","key":314,"tp1":"cu","vl1":"EN","cn1":80}'
    7 = ',{"id":"ev5_oneplustwo_2","key":537,"tp1":"cu","vl1":"EN","cn1":300},{"id":"ev5_oneplustwo_3","key":314,"tp1":"cu","vl1":"EN","cn1":80}'
    8 = ',{"id":"ev5_oneplustwo_2","key":537,"tp1":"cu","vl1":"EN","cn1":300},{"id":"ev5_oneplustwo_3","key":314,"tp1":"cu","vl1":"EN","cn1":80}'
    9 = ',{"id":"ev5_oneplustwo_2","key":537,"tp1":"cu","vl1":"EN","cn1":300},{"id":"ev5_oneplustwo_3","key":314,"tp1":"cu","vl1":"EN","cn1":80}'
    10 = ',{"id":"ev5_oneplustwo_2","key":537,"tp1":"cu","vl1":"EN","cn1":300},{"id":"ev5_oneplustwo_3","key":314,"tp1":"cu","vl1":"EN","cn1":80},{"id":"ev4_conti_1","key":721,"tp1":"cu","vl1":"EN","cn1":80,"tp2":"cu","vl2":"GO","cn2":35000,"tp3":"cu","vl3":"EN","cn3":170}'
    11 = ',{"id":"ev5_oneplustwo_2","key":537,"tp1":"cu","vl1":"EN","cn1":300},{"id":"ev5_oneplustwo_3","key":314,"tp1":"cu","vl1":"EN","cn1":80},{"id":"ev4_conti_1","key":721,"tp1":"cu","vl1":"EN","cn1":80,"tp2":"cu","vl2":"GO","cn2":35000,"tp3":"cu","vl3":"EN","cn3":170},{"id":"ev4_conti_2","key":884,"tp1":"cu","vl1":"EN","cn1":150}'
    12 = ',{"id":"ev5_oneplustwo_2","key":537,"tp1":"cu","vl1":"EN","cn1":300},{"id":"ev5_oneplustwo_3","key":314,"tp1":"cu","vl1":"EN","cn1":80},{"id":"ev4_conti_1","key":721,"tp1":"cu","vl1":"EN","cn1":80,"tp2":"cu","vl2":"GO","cn2":35000,"tp3":"cu","vl3":"EN","cn3":170},{"id":"ev4_conti_2","key":884,"tp1":"cu","vl1":"EN","cn1":150}'
    13 = ',{"id":"ev5_oneplustwo_2","key":537,"tp1":"cu","vl1":"EN","cn1":300},{"id":"ev5_oneplustwo_3","key":314,"tp1":"cu","vl1":"EN","cn1":80},{"id":"ev4_conti_1","key":721,"tp1":"cu","vl1":"EN","cn1":80,"tp2":"cu","vl2":"GO","cn2":35000,"tp3":"cu","vl3":"EN","cn3":170},{"id":"ev4_conti_2","key":884,"tp1":"cu","vl1":"EN","cn1":150},{"id":"ev4_conti_4","key":394,"tp1":"cu","vl1":"EN","cn1":150,"tp2":"cu","vl2":"GO","cn2":20000}'
    14 = ',{"id":"ev5_oneplustwo_2","key":537,"tp1":"cu","vl1":"EN","cn1":300},{"id":"ev5_oneplustwo_3","key":314,"tp1":"cu","vl1":"EN","cn1":80},{"id":"ev4_conti_1","key":721,"tp1":"cu","vl1":"EN","cn1":80,"tp2":"cu","vl2":"GO","cn2":35000,"tp3":"cu","vl3":"EN","cn3":170},{"id":"ev4_conti_2","key":884,"tp1":"cu","vl1":"EN","cn1":150},{"id":"ev4_conti_4","key":394,"tp1":"cu","vl1":"EN","cn1":150,"tp2":"cu","vl2":"GO","cn2":20000},{"id":"ev4_conti_5","key":612,"tp1":"cu","vl1":"EN","cn1":80,"tp2":"cu","vl2":"GO","cn2":10000,"tp3":"cu","vl3":"EN","cn3":200}'
    15 = ',{"id":"ev5_oneplustwo_2","key":537,"tp1":"cu","vl1":"EN","cn1":300},{"id":"ev5_oneplustwo_3","key":314,"tp1":"cu","vl1":"EN","cn1":80},{"id":"ev4_conti_1","key":721,"tp1":"cu","vl1":"EN","cn1":80,"tp2":"cu","vl2":"GO","cn2":35000,"tp3":"cu","vl3":"EN","cn3":170},{"id":"ev4_conti_2","key":884,"tp1":"cu","vl1":"EN","cn1":150},{"id":"ev4_conti_4","key":394,"tp1":"cu","vl1":"EN","cn1":150,"tp2":"cu","vl2":"GO","cn2":20000},{"id":"ev4_conti_5","key":612,"tp1":"cu","vl1":"EN","cn1":80,"tp2":"cu","vl2":"GO","cn2":10000,"tp3":"cu","vl3":"EN","cn3":200}'
    16 = ',{"id":"ev5_oneplustwo_2","key":537,"tp1":"cu","vl1":"EN","cn1":300},{"id":"ev5_oneplustwo_3","key":314,"tp1":"cu","vl1":"EN","cn1":80},{"id":"ev4_conti_1","key":721,"tp1":"cu","vl1":"EN","cn1":80,"tp2":"cu","vl2":"GO","cn2":35000,"tp3":"cu","vl3":"EN","cn3":170},{"id":"ev4_conti_2","key":884,"tp1":"cu","vl1":"EN","cn1":150},{"id":"ev4_conti_4","key":394,"tp1":"cu","vl1":"EN","cn1":150,"tp2":"cu","vl2":"GO","cn2":20000},{"id":"ev4_conti_5","key":612,"tp1":"cu","vl1":"EN","cn1":80,"tp2":"cu","vl2":"GO","cn2":10000,"tp3":"cu","vl3":"EN","cn3":200},{"id":"ev4_conti_7","key":930,"tp1":"cu","vl1":"GO","cn1":50000}'
    17 = ',{"id":"ev5_oneplustwo_2","key":537,"tp1":"cu","vl1":"EN","cn1":300},{"id":"ev5_oneplustwo_3","key":314,"tp1":"cu","vl1":"EN","cn1":80},{"id":"ev4_conti_1","key":721,"tp1":"cu","vl1":"EN","cn1":80,"tp2":"cu","vl2":"GO","cn2":35000,"tp3":"cu","vl3":"EN","cn3":170},{"id":"ev4_conti_2","key":884,"tp1":"cu","vl1":"EN","cn1":150},{"id":"ev4_conti_4","key":394,"tp1":"cu","vl1":"EN","cn1":150,"tp2":"cu","vl2":"GO","cn2":20000},{"id":"ev4_conti_5","key":612,"tp1":"cu","vl1":"EN","cn1":80,"tp2":"cu","vl2":"GO","cn2":10000,"tp3":"cu","vl3":"EN","cn3":200},{"id":"ev4_conti_7","key":930,"tp1":"cu","vl1":"GO","cn1":50000},{"id":"ev4_conti_8","key":959,"tp1":"cu","vl1":"EN","cn1":350,"tp2":"cu","vl2":"GO","cn2":30000}'
}

$xVals = @{
    5 = '{"id":"ev5_oneplustwo_2","key":537,"tp1":"cu","vl1":"EN","cn1":300}'
    6 = '{"id":"ev5_oneplustwo_3","key":314,"tp1":"cu","vl1":"EN","cn1":80}'
    10 = '{"id":"ev4_conti_1","key":721,"tp1":"cu","vl1":"EN","cn1":80,"tp2":"cu","vl2":"GO","cn2":35000,"tp3":"cu","vl3":"EN","cn3":170}'
    11 = '{"id":"ev4_conti_2","key":884,"tp1":"cu","vl1":"EN","cn1":150}'
    13 = '{"id":"ev4_conti_4","key":394,"tp1":"cu","vl1":"EN","cn1":150,"tp2":"cu","vl2":"GO","cn2":20000}'
    14 = '{"id":"ev4_conti_5","key":612,"tp1":"cu","vl1":"EN","cn1":80,"tp2":"cu","vl2":"GO","cn2":10000,"tp3":"cu","vl3":"EN","cn3":200}'
    16 = '{"id":"ev4_conti_7","key":930,"tp1":"cu","vl1":"GO","cn1":50000}'
    17 = '{"id":"ev4_conti_8","key":959,"tp1":"cu","vl1":"EN","cn1":350,"tp2":"cu","vl2":"GO","cn2":30000}'
}

foreach ($r in $wVals.Keys) {
    $ws.Cells.Item($r, 23).Value = $wVals[$r]
}
foreach ($r in $xVals.Keys) {
    $ws.Cells.Item($r, 24).Value = $xVals[$r]
}

# --- 3. Reset the saved selection to A1 (matches the post-edit save state) ---
$ws.Range("A1").Select()
